$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

$ws.Range("A107").Value = 5
$ws.Range("B107").Value = "Macroferia Regional de Talca"
$ws.Range("C107").Value = "Maule"
$ws.Range("D107").Value = 44873
$ws.Range("E107").Value = 7
$ws.Range("F107").Value = 100112028
$ws.Range("G107").Value = "Sandia"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 750
$ws.Range("N107").Value = "$/kilo"
$ws.Range("O107").Value = "Perú"
$ws.Range("P107").Value = 750
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = "Hortaliza"
